# "Generate Report for Handback" -- refresh the localization-status report:
#  - Overview sheet: status flips from "Ready for handoff" to
#    "Handed back: in sync with en-US" for both locales.
#  - zh-cn / de-de sheets: the handback timestamp advances and the stale
#    "handback file is not latest" error is cleared now that it's in sync.
#  - A few report columns are widened/narrowed to fit the refreshed content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview: handback status for both locales ---
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn: refreshed handback datetime + cleared error detail ---
$zhcn.Range("K2").Value = "2016-08-22 00:59:38"
$zhcn.Range("P2").Value = ""

# --- de-de: refreshed handback datetime + cleared error detail ---
$dede.Range("K2").Value = "2016-08-22 00:59:44"
$dede.Range("P2").Value = ""

# --- Column width adjustments to fit the refreshed report text ---
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
